$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Paragraph 1 ("Cost and service call to action "): strike the whole
#    paragraph (adds <w:strike/> to the paragraph mark rPr and to every run's
#    rPr, creating a new rPr for the first run which previously had none).
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$p1.Range.Font.StrikeThrough = 1

# ---------------------------------------------------------------------------
# 2) Paragraph 2 ("Lab devices big pics are wrong models "): append a new
#    run "--- review this " and place a collapsed _GoBack bookmark right
#    after it (moved from what used to be the following empty paragraph).
#
#    NOTE: the COM host mis-places a bookmark that is Add()-ed as a
#    *collapsed* range sitting exactly on a paragraph-mark boundary (it can
#    snap to an unrelated paragraph). To avoid that, we briefly insert a
#    placeholder character, wrap it with the bookmark (a normal, non-
#    collapsed, non-boundary range - which works correctly), and then delete
#    the placeholder; deleting the wrapped text naturally collapses the
#    bookmark to the correct spot, exactly as real Word does.
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(2)
$r2 = $p2.Range
$ip2 = $d.Range($r2.End - 1, $r2.End - 1)
$ip2.InsertAfter("--- review this X")

$p2b = $d.Paragraphs.Item(2)
$rb2 = $p2b.Range
$phStart = $rb2.End - 2
$phEnd = $rb2.End - 1
$phRange = $d.Range($phStart, $phEnd)
$d.Bookmarks.Add("_GoBack", $phRange)
$placeholder = $d.Range($phStart, $phEnd)
$placeholder.Delete()

# ---------------------------------------------------------------------------
# 3) Paragraph 3 ("Bigger font lab devices bigger pics"): append a new bold
#    run " --- review this ".
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
$r3 = $p3.Range
$ip3 = $d.Range($r3.End - 1, $r3.End - 1)
$ip3.InsertAfter(" --- review this ")
$ip3.Font.Bold = 1
$ip3.Font.BoldBi = 1

Write-Output "done"
